$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 492; this pushes the existing rows
# 492-545 down to 493-546 (dimension grows from A1:R545 to A1:R546).
$ws.Rows.Item(492).Insert()

# Populate the newly inserted row 492 with the new price-report entry.
$ws.Range("A492").Value = 5
$ws.Range("B492").Value = "Macroferia Regional de Talca"
$ws.Range("C492").Value = "Maule"
$ws.Range("D492").Value = 45194
$ws.Range("E492").Value = 7
$ws.Range("F492").Value = 100112003
$ws.Range("G492").Value = "Ajo"
$ws.Range("H492").Value = "Chino"
$ws.Range("I492").Value = "Primera"
$ws.Range("J492").Value = 300
$ws.Range("K492").Value = 21000
$ws.Range("L492").Value = 21000
$ws.Range("M492").Value = 21000
$ws.Range("N492").Value = "`$/malla 10 kilos"
$ws.Range("O492").Value = "China"
$ws.Range("P492").Value = 2100
$ws.Range("Q492").Value = 10
$ws.Range("R492").Value = "Hortaliza"
